$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: append a red parenthetical note to the first paragraph.
# "This is a Microsoft word document." ->
#   "This is a Microsoft word document.  " (plain, 2 trailing spaces)
#   + "(This is a change – Ve" (red)
#   + "rsion for main branch" (red)
#   + ")" (red)
# ------------------------------------------------------------------
$p = $d.Paragraphs(1)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)

$r.InsertAfter("  ")
$r.Collapse(0)

$r.InsertAfter("(This is a change – Ve")
$r.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter("rsion for main branch")
$r.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter(")")
$r.Font.Color = 255
$r.Collapse(0)

# ------------------------------------------------------------------
# Edit 2: drop the trailing "ank God almighty, we are free at last."
# paragraph (the tail end of the Raven poem section).
# ------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.Delete()

# ------------------------------------------------------------------
# Edit 3: prune the unused styles that Word drops on save once
# nothing in the body references them any more.
# ------------------------------------------------------------------
$unusedStyles = @(
  "Heading 2",
  "Heading 4",
  "apple-converted-space",
  "Hyperlink",
  "Heading 2 Char",
  "Heading 4 Char",
  "audio-tool",
  "subscribe",
  "subscribe-more-info",
  "generic-title",
  "podcast-tools__subscribe-links"
)
for ($i = $d.Styles.Count; $i -ge 1; $i--) {
  $s = $d.Styles.Item($i)
  if ($unusedStyles -contains $s.NameLocal) {
    $s.Delete()
  }
}

Write-Output ("Para1: " + $d.Paragraphs(1).Range.Text)
Write-Output ("ParaCount: " + $d.Paragraphs.Count)
Write-Output ("LastPara: " + $d.Paragraphs($d.Paragraphs.Count).Range.Text)
Write-Output ("StyleCount: " + $d.Styles.Count)
